$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("N2").Value = 70831.95557958097
$ws2025.Range("O2").Value = 69610.44223910036

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 47338.61932520662
$ws2030.Range("I2").Value = 283167.7766510943
$ws2030.Range("L2").Value = 178095.3756971828
$ws2030.Range("M2").Value = 114008.3253427963
$ws2030.Range("N2").Value = 33931.82461160053
$ws2030.Range("O2").Value = 50657.26889981552

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 29273.60317916481
$ws2035.Range("B2").Value = 22330.72247668595
$ws2035.Range("E2").Value = 110739.3594843864
$ws2035.Range("I2").Value = 150386.9441391908
$ws2035.Range("M2").Value = 35556.98862372932
$ws2035.Range("N2").Value = 44813.41193308897
$ws2035.Range("O2").Value = 26775.55841092002

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("N2").Value = 1041.156112142704

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 34256.25291363284
$ws2045.Range("N2").Value = 5271.89502409355
$ws2045.Range("O2").Value = 22972.54525065989
